# Fruta / hortaliza, semanal
# Insert a new weekly record as row 15 in the price history table, shifting
# all subsequent rows (formerly 15-41) down by one to (16-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15; Excel pushes rows 15..41 down to 16..42
# and extends the used range (dimension becomes A1:R42) automatically.
$ws.Rows("15:15").Insert()

# Populate the freshly inserted row 15. Columns A, B, C, E, F, G, H, I, N, O,
# Q, R repeat the same market/category metadata as every other row in this
# subset, and K/L (min/max price) stay at 14000 / 15000 for this new record.
$ws.Range("A15").Value2 = 3
$ws.Range("B15").Value2 = "Femacal de La Calera"
$ws.Range("C15").Value2 = "Coquimbo"
$ws.Range("D15").Value2 = 44725
$ws.Range("E15").Value2 = 5
$ws.Range("F15").Value2 = 100112035
$ws.Range("G15").Value2 = "Bruselas (repollito)"
$ws.Range("H15").Value2 = "Sin especificar"
$ws.Range("I15").Value2 = "Primera"
$ws.Range("J15").Value2 = 85
$ws.Range("K15").Value2 = 14000
$ws.Range("L15").Value2 = 15000
$ws.Range("M15").Value2 = 14471
$ws.Range("N15").Value2 = "$/malla 15 kilos"
$ws.Range("O15").Value2 = "Provincia de Quillota"
$ws.Range("P15").Value2 = 965
$ws.Range("Q15").Value2 = 15
$ws.Range("R15").Value2 = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
